$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D are written as text,
# matching the original inline-string cell type (prevents Excel auto-
# converting values like "7.10" or "2.00" into numbers and losing precision).
$priceCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D21","D23","D24","D27","D28","D29","D30","D31","D32","D33","D34","D35","D37","D38","D40","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "47.250.20"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.489.44"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D5").Value = "321.09"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "107.75"
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("D7").Value = "0.521"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").Value = "38.52"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "18.33"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "7.10"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "2.870.77"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "2.494.06"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "0.847"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "47.146.42"
$ws.Range("D19").Value = "12.78"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "6.60"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "0.0₃0930"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("E22").Value = "  +14.29%  "
$ws.Range("D23").Value = "70.19"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D24").Value = "245.18"
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "25.67"
$ws.Range("E27").Value = "  -3.37%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "34.48"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "0.135"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").Value = "49.48"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "20.38"
$ws.Range("E33").Value = "  +2.62%  "
$ws.Range("D34").Value = "5.31"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").Value = "0.0778"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D37").Value = "1.96"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "4.62"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").Value = "22.70"
$ws.Range("E40").Value = "  +4.44%  "
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "118.33"
$ws.Range("E43").Value = "  -3.99%  "
$ws.Range("D44").Value = "0.0295"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = "1.982.89"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "3.01"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").Value = "2.00"
$ws.Range("E47").Value = "  -6.28%  "
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("D49").Value = "1.77"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").Value = "5.09"
$ws.Range("E50").Value = "  -6.97%  "
$ws.Range("D51").Value = "56.50"
$ws.Range("E51").Value = "  +2.83%  "
